$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text (prices, percentages) is stored as literal text,
# matching the source inlineStr cells (avoid Excel auto-converting to number/percent).
$textCells = @("D2","E2","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E19","D20","E20","E21","D22","E22","D23","E23","E24","D25","E25","D26","E26","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "310.07"
$ws.Range("E2").Value = "-0.23%"

$ws.Range("E3").Value = "-0.13%"

$ws.Range("D4").Value = "5.184"
$ws.Range("E4").Value = "1.60%"

$ws.Range("D5").Value = "0.07871"
$ws.Range("E5").Value = "1.53%"

$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "4.423"
$ws.Range("E6").Value = "1.76%"

$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "1.911"
$ws.Range("E7").Value = "1.85%"

$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "8.262"
$ws.Range("E8").Value = "0.72%"

$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "3.001"
$ws.Range("E9").Value = "4.14%"

$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "0.9369"
$ws.Range("E10").Value = "1.85%"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "0.1120"
$ws.Range("E11").Value = "-5.94%"

$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "0.1971"
$ws.Range("E12").Value = "2.97%"

$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "0.09102"
$ws.Range("E13").Value = "2.37%"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03289"
$ws.Range("E14").Value = "-2.81%"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09603"
$ws.Range("E15").Value = "-0.84%"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001378"
$ws.Range("E16").Value = "0.75%"

$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.006114"
$ws.Range("E17").Value = "5.74%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.613"
$ws.Range("E18").Value = "1.65%"

$ws.Range("E19").Value = "0.12%"

$ws.Range("D20").Value = "6.457"
$ws.Range("E20").Value = "28.22%"

$ws.Range("E21").Value = "1.63%"

$ws.Range("D22").Value = "0.2518"
$ws.Range("E22").Value = "-2.79%"

$ws.Range("D23").Value = "0.04391"
$ws.Range("E23").Value = "0.08%"

$ws.Range("E24").Value = "1.77%"

$ws.Range("D25").Value = "0.004580"
$ws.Range("E25").Value = "7.88%"

$ws.Range("D26").Value = "0.0001362"
$ws.Range("E26").Value = "0.84%"

$ws.Range("D39").Value = "0.02211"
$ws.Range("E39").Value = "5.58%"

$ws.Range("D40").Value = "0.05104"
$ws.Range("E40").Value = "3.21%"

$ws.Range("D41").Value = "0.007462"
$ws.Range("E41").Value = "-2.61%"

$ws.Range("D42").Value = "0.1354"
$ws.Range("E42").Value = "1.05%"

$ws.Range("D43").Value = "0.008759"
$ws.Range("E43").Value = "-11.29%"

$ws.Range("D44").Value = "0.002133"
$ws.Range("E44").Value = "3.50%"

$ws.Range("D45").Value = "0.008627"
$ws.Range("E45").Value = "-10.35%"

$ws.Range("D46").Value = "0.00006562"
$ws.Range("E46").Value = "-0.06%"

$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.15%"

$ws.Range("D48").Value = "0.002862"
$ws.Range("E48").Value = "-5.94%"

$ws.Range("D49").Value = "0.001002"
$ws.Range("E49").Value = "-40.68%"

$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "0.15%"

$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "0.15%"
